$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price column cells being updated so Excel
# does not reinterpret numeric-looking strings (e.g. "5.46", "59.268.32")
# as numbers, then restore the default (no explicit) style afterwards so the
# cell formatting matches the original inline-string cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.268.32"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.521.28"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.89"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.97"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  -2.95%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("E8").Value = "  -1.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.526.73"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("E11").Value = "  +1.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.46"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  +1.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.967.14"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.193.78"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.89"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  -2.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.523.04"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.94"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.24"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.38"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.20"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("E25").Value = "  -2.83%  "

$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.79"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0767"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.59"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +2.85%  "

$ws.Range("E33").Value = "  +0.29%  "

$ws.Range("E34").Value = "  +1.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.13"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  -6.41%  "

$ws.Range("E36").Value = "  -0.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.22"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  -3.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.58"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  -2.49%  "

$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "284.05"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  -5.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.25"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  -6.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.86"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.597"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0930"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.64"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  -2.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.52"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0510"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -0.89%  "

$ws.Range("E51").Value = "  -1.87%  "
